# Applies the FFXIV leve-profit data refresh across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner update.
# Only numeric H:N (price/profit) columns are touched; H:L are refreshed
# market-board prices, M/N are the derived NQ/HQ profit figures (stored as
# plain values, not formulas, in this workbook).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2524.4
$ws.Range("I6").Value = 2524.4
$ws.Range("K6").Value = 7573.200000000001
$ws.Range("M6").Value = -7461.200000000001
$ws.Range("H11").Value = 76.69231
$ws.Range("I11").Value = 76.69231
$ws.Range("K11").Value = 76.69231
$ws.Range("M11").Value = 63.30768999999999
$ws.Range("H17").Value = 975
$ws.Range("I17").Value = 960
$ws.Range("J17").Value = 990
$ws.Range("K17").Value = 2880
$ws.Range("L17").Value = 2970
$ws.Range("M17").Value = -2712
$ws.Range("N17").Value = -3306
$ws.Range("H127").Value = 5858.5
$ws.Range("I127").Value = 6706.0625
$ws.Range("J127").Value = 3598.3333
$ws.Range("K127").Value = 20118.1875
$ws.Range("L127").Value = 10794.9999
$ws.Range("M127").Value = -15158.1875
$ws.Range("N127").Value = -20714.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2234729.2
$ws.Range("I32").Value = 2787210
$ws.Range("J32").Value = 24805.555
$ws.Range("K32").Value = 2787210
$ws.Range("L32").Value = 24805.555
$ws.Range("M32").Value = -2786923
$ws.Range("N32").Value = -25379.555
$ws.Range("H62").Value = 50249
$ws.Range("J62").Value = 50249
$ws.Range("L62").Value = 50249
$ws.Range("N62").Value = -51497
$ws.Range("H63").Value = 1883.2667
$ws.Range("I63").Value = 1875
$ws.Range("K63").Value = 1875
$ws.Range("M63").Value = -1189
$ws.Range("H65").Value = 50249
$ws.Range("J65").Value = 50249
$ws.Range("L65").Value = 150747
$ws.Range("N65").Value = -156987
$ws.Range("H66").Value = 1883.2667
$ws.Range("I66").Value = 1875
$ws.Range("K66").Value = 9375
$ws.Range("M66").Value = -5943
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H132").Value = 25653554
$ws.Range("I132").Value = 27618216
$ws.Range("J132").Value = 11115051
$ws.Range("K132").Value = 82854648
$ws.Range("L132").Value = 33345153
$ws.Range("M132").Value = -82852118
$ws.Range("N132").Value = -33350213

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 978.8
$ws.Range("I11").Value = 978.8
$ws.Range("K11").Value = 978.8
$ws.Range("M11").Value = -838.8
$ws.Range("H54").Value = 4482.6
$ws.Range("I54").Value = 603.25
$ws.Range("J54").Value = 20000
$ws.Range("K54").Value = 603.25
$ws.Range("L54").Value = 20000
$ws.Range("M54").Value = -119.25
$ws.Range("N54").Value = -20968

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 170.5
$ws.Range("I7").Value = 300
$ws.Range("J7").Value = 127.333336
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 127.333336
$ws.Range("M7").Value = -187
$ws.Range("N7").Value = -353.333336
$ws.Range("H59").Value = 16450.375
$ws.Range("J59").Value = 16450.375
$ws.Range("L59").Value = 16450.375
$ws.Range("N59").Value = -18740.375
$ws.Range("H134").Value = 1214620.2
$ws.Range("I134").Value = 1748.7
$ws.Range("J134").Value = 3080576.5
$ws.Range("K134").Value = 5246.1
$ws.Range("L134").Value = 9241729.5
$ws.Range("M134").Value = -2711.1
$ws.Range("N134").Value = -9246799.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5274649
$ws.Range("J131").Value = 1054.8096
$ws.Range("L131").Value = 3164.4288
$ws.Range("N131").Value = -13244.4288

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 32000
$ws.Range("J64").Value = 32000
$ws.Range("L64").Value = 32000
$ws.Range("N64").Value = -32496
$ws.Range("H67").Value = 32000
$ws.Range("J67").Value = 32000
$ws.Range("L67").Value = 32000
$ws.Range("N67").Value = -33716
$ws.Range("H132").Value = 8061777
$ws.Range("I132").Value = 7739427.5
$ws.Range("K132").Value = 23218282.5
$ws.Range("M132").Value = -23215752.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 35717584
$ws.Range("I22").Value = 950
$ws.Range("J22").Value = 41670360
$ws.Range("K22").Value = 950
$ws.Range("L22").Value = 41670360
$ws.Range("M22").Value = -655
$ws.Range("N22").Value = -41670950
$ws.Range("H27").Value = 35717584
$ws.Range("I27").Value = 950
$ws.Range("J27").Value = 41670360
$ws.Range("K27").Value = 950
$ws.Range("L27").Value = 41670360
$ws.Range("M27").Value = -843
$ws.Range("N27").Value = -41670574
$ws.Range("H42").Value = 16009.333
$ws.Range("J42").Value = 16009.333
$ws.Range("L42").Value = 16009.333
$ws.Range("N42").Value = -17135.333
$ws.Range("H46").Value = 62501036
$ws.Range("I46").Value = 899.75
$ws.Range("J46").Value = 125001176
$ws.Range("K46").Value = 899.75
$ws.Range("L46").Value = 125001176
$ws.Range("M46").Value = -711.75
$ws.Range("N46").Value = -125001552
$ws.Range("H49").Value = 16009.333
$ws.Range("J49").Value = 16009.333
$ws.Range("L49").Value = 16009.333
$ws.Range("N49").Value = -16303.333
$ws.Range("H55").Value = 10869659
$ws.Range("I55").Value = 16666732
$ws.Range("J55").Value = 145.625
$ws.Range("K55").Value = 16666732
$ws.Range("L55").Value = 145.625
$ws.Range("M55").Value = -16666559
$ws.Range("N55").Value = -491.625
$ws.Range("H64").Value = 20000
$ws.Range("J64").Value = 20000
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20450
$ws.Range("H67").Value = 20000
$ws.Range("J67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
